$wb = $excel.ActiveWorkbook

# Rows (sheet row numbers) and their updated "F" column (想去人数) values.
$updates = @{
    4  = 1531
    5  = 585
    6  = 1077
    7  = 11181
    10 = 330
    12 = 769
    13 = 12263
    14 = 12856
    15 = 33
    21 = 59
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
